$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new titration row (20220314) to the CRM accuracy log ---
$ws.Range("A66").Value = 20220314
$ws.Range("B66").Value = 2221.5278341974199
$ws.Range("C66").Value = 2224.4699999999998
$ws.Range("D66").Formula = "=100*(B66-C66)/C66"
$ws.Range("E66").Value = 180
$ws.Range("F66").Value = "CRM OPENED 20220302"

# --- Give the new "Batch value" column a bit more breathing room ---
$ws.Columns.Item(2).ColumnWidth = 14.29

# --- Scroll the sheet to show the newly-added rows / selection ---
$excel.Goto($ws.Range("A56"), $true) | Out-Null
$ws.Range("E69").Select() | Out-Null
